$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column E (Param3) entirely - shifts dimension to A1:D26
$ws.Range("E:E").Delete()

# Update data rows 2-26 with new file names, class names, and Param1 values
$ws.Range("A2").Value = 'Zelda--param1-00.99--5-01.dac'
$ws.Range("B2").Value = 'Zelda'
$ws.Range("C2").Value = 0.99
$ws.Range("D2").Value = 1

$ws.Range("A3").Value = 'Zelda--param1-00.63--6-08.dac'
$ws.Range("B3").Value = 'Zelda'
$ws.Range("C3").Value = 0.63
$ws.Range("D3").Value = 8

$ws.Range("A4").Value = 'Zelda--param1-00.67--5-20.dac'
$ws.Range("B4").Value = 'Zelda'
$ws.Range("C4").Value = 0.67
$ws.Range("D4").Value = 20

$ws.Range("A5").Value = 'Zelda--param1-00.85--4-20.dac'
$ws.Range("B5").Value = 'Zelda'
$ws.Range("C5").Value = 0.85
$ws.Range("D5").Value = 20

$ws.Range("A6").Value = 'Zelda--param1-00.89--1-35.dac'
$ws.Range("B6").Value = 'Zelda'
$ws.Range("C6").Value = 0.89
$ws.Range("D6").Value = 35

$ws.Range("A7").Value = 'lofi--param1-00.58.dac'
$ws.Range("B7").Value = 'lofi'
$ws.Range("C7").Value = 0.58

$ws.Range("A8").Value = 'lofi--param1-00.49.dac'
$ws.Range("B8").Value = 'lofi'
$ws.Range("C8").Value = 0.49

$ws.Range("A9").Value = 'lofi--param1-00.88.dac'
$ws.Range("B9").Value = 'lofi'
$ws.Range("C9").Value = 0.88

$ws.Range("A10").Value = 'lofi--param1-00.44.dac'
$ws.Range("B10").Value = 'lofi'
$ws.Range("C10").Value = 0.44

$ws.Range("A11").Value = 'lofi--param1-00.81.dac'
$ws.Range("B11").Value = 'lofi'
$ws.Range("C11").Value = 0.81

$ws.Range("A12").Value = 'Fusion--param1-00.41.dac'
$ws.Range("B12").Value = 'Fusion'
$ws.Range("C12").Value = 0.41

$ws.Range("A13").Value = 'Fusion--param1-00.91.dac'
$ws.Range("B13").Value = 'Fusion'
$ws.Range("C13").Value = 0.91

$ws.Range("A14").Value = 'Fusion--param1-00.76.dac'
$ws.Range("B14").Value = 'Fusion'
$ws.Range("C14").Value = 0.76

$ws.Range("A15").Value = 'Fusion--param1-00.19.dac'
$ws.Range("B15").Value = 'Fusion'
$ws.Range("C15").Value = 0.19

$ws.Range("A16").Value = 'Fusion--param1-00.67.dac'
$ws.Range("B16").Value = 'Fusion'
$ws.Range("C16").Value = 0.67

$ws.Range("A17").Value = '8bit--param1-00.22.dac'
$ws.Range("B17").Value = '8bit'
$ws.Range("C17").Value = 0.22

$ws.Range("A18").Value = '8bit--param1-00.71--04.dac'
$ws.Range("B18").Value = '8bit'
$ws.Range("C18").Value = 0.71

$ws.Range("A19").Value = '8bit--param1-00.67--09.dac'
$ws.Range("B19").Value = '8bit'
$ws.Range("C19").Value = 0.67

$ws.Range("A20").Value = '8bit--param1-00.44.dac'
$ws.Range("B20").Value = '8bit'
$ws.Range("C20").Value = 0.44

$ws.Range("A21").Value = '8bit--param1-00.54.dac'
$ws.Range("B21").Value = '8bit'
$ws.Range("C21").Value = 0.54

$ws.Range("A22").Value = 'duduk--param1-00.67.dac'
$ws.Range("B22").Value = 'duduk'
$ws.Range("C22").Value = 0.67

$ws.Range("A23").Value = 'duduk--param1-00.03.dac'
$ws.Range("B23").Value = 'duduk'
$ws.Range("C23").Value = 0.03

$ws.Range("A24").Value = 'duduk--param1-00.54--KSHMR_Duduk_19_One_Shot_F#m.dac'
$ws.Range("B24").Value = 'duduk'
$ws.Range("C24").Value = 0.54

$ws.Range("A25").Value = 'duduk--param1-00.05.dac'
$ws.Range("B25").Value = 'duduk'
$ws.Range("C25").Value = 0.05

$ws.Range("A26").Value = 'duduk--param1-00.30.dac'
$ws.Range("B26").Value = 'duduk'
$ws.Range("C26").Value = 0.3
